$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("G7").Value = 1.95
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 4.33
$ws.Range("J7").Value = 2.75
$ws.Range("K7").Value = 1.95
$ws.Range("L7").Value = 5
$ws.Range("Q7").Value = 2.5
$ws.Range("R7").Value = 1.5
$ws.Range("S7").Value = 4.2
$ws.Range("T7").Value = 1.23
$ws.Range("U7").Value = 5
$ws.Range("V7").Value = 1.17
$ws.Range("AB7").Value = 7.5
$ws.Range("AC7").Value = 9.5
$ws.Range("AD7").Value = 17
$ws.Range("AE7").Value = 19
$ws.Range("AH7").Value = 6.5
$ws.Range("AI7").Value = 21
$ws.Range("AL7").Value = 9
$ws.Range("AM7").Value = 21
$ws.Range("AN7").Value = 15
$ws.Range("AO7").Value = 51
$ws.Range("AP7").Value = 41
$ws.Range("AR7").Value = 1.93
$ws.Range("AS7").Value = 1.93

# Row 8
$ws.Range("G8").Value = 1.8
$ws.Range("I8").Value = 4.75
$ws.Range("K8").Value = 1.95
$ws.Range("L8").Value = 5.5
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 1.5
$ws.Range("P8").Value = 2.5
$ws.Range("Q8").Value = 2.5
$ws.Range("R8").Value = 1.5
$ws.Range("S8").Value = 4.1
$ws.Range("T8").Value = 1.24
$ws.Range("W8").Value = 1.57
$ws.Range("X8").Value = 2.25
$ws.Range("AE8").Value = 19
$ws.Range("AH8").Value = 6.5
$ws.Range("AL8").Value = 9.5
$ws.Range("AP8").Value = 41
$ws.Range("AR8").Value = 1.9
$ws.Range("AS8").Value = 1.95

# Row 9
$ws.Range("Y9").Value = 2.1
$ws.Range("Z9").Value = 1.67
$ws.Range("AO9").Value = 29

# Row 11
$ws.Range("M11").Value = 1.17
$ws.Range("N11").Value = 5
$ws.Range("W11").Value = 1.75
$ws.Range("X11").Value = 2.05

# Row 12
$ws.Range("H12").Value = 2.88
$ws.Range("I12").Value = 6.25
$ws.Range("J12").Value = 2.6
$ws.Range("L12").Value = 7
$ws.Range("N12").Value = 4.75
$ws.Range("AA12").Value = 4.33
$ws.Range("AI12").Value = 29
$ws.Range("AJ12").Value = 151
$ws.Range("AM12").Value = 29
$ws.Range("AO12").Value = 81
$ws.Range("AP12").Value = 67

# Row 24
$ws.Range("G24").Value = 3.2
$ws.Range("I24").Value = 2.3
$ws.Range("J24").Value = 3.75
$ws.Range("L24").Value = 3
$ws.Range("AC24").Value = 12
$ws.Range("AD24").Value = 34
$ws.Range("AO24").Value = 21
$ws.Range("AP24").Value = 19
$ws.Range("AQ24").Value = 29

# Row 27
$ws.Range("G27").Value = 3.7
$ws.Range("I27").Value = 2.1
$ws.Range("Y27").Value = 2.2
$ws.Range("Z27").Value = 1.62
$ws.Range("AC27").Value = 13
$ws.Range("AM27").Value = 8.5
$ws.Range("AN27").Value = 10
$ws.Range("AO27").Value = 19

# Row 28
$ws.Range("G28").Value = 3.1
$ws.Range("K28").Value = 1.8
$ws.Range("M28").Value = 1.13
$ws.Range("N28").Value = 6
$ws.Range("S28").Value = 5
$ws.Range("T28").Value = 1.16
$ws.Range("U28").Value = 6
$ws.Range("V28").Value = 1.13
$ws.Range("AI28").Value = 23
$ws.Range("AP28").Value = 29

# Row 36
$ws.Range("G36").Value = 2.25
$ws.Range("I36").Value = 3.2
$ws.Range("AE36").Value = 21
$ws.Range("AK36").Value = 301

# Row 70
$ws.Range("M70").Value = 1.07
$ws.Range("N70").Value = 9
$ws.Range("U70").Value = 4
$ws.Range("V70").Value = 1.22

# Row 71
$ws.Range("G71").Value = 1.5
$ws.Range("I71").Value = 7
$ws.Range("L71").Value = 7
$ws.Range("O71").Value = 1.33
$ws.Range("P71").Value = 3.25
$ws.Range("Q71").Value = 2.05
$ws.Range("R71").Value = 1.75
$ws.Range("Y71").Value = 2.2
$ws.Range("Z71").Value = 1.62
$ws.Range("AC71").Value = 9
$ws.Range("AD71").Value = 10
$ws.Range("AE71").Value = 13
$ws.Range("AG71").Value = 8.5
$ws.Range("AL71").Value = 15
$ws.Range("AM71").Value = 34
$ws.Range("AO71").Value = 81
